# Neutralize racial-voter language in the non-electoral resume.
# Three locations change "affecting all Black and Asian-American voters"
# to "affecting 50M voters" (plain text in two spots, with a bold/colored
# "50M" run matching the existing numeric-highlight style in the bullet).

$d = $word.ActiveDocument

# --- Edit 1: Professional summary paragraph (plain text run) ---
$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed",
    2) | Out-Null

# --- Edit 3: Project "Impact:" statement (plain text run) ---
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved",
    2) | Out-Null

# --- Edit 2: Experience bullet - needs "50M" split into its own
#     bold/colored run (matching the "23%"/"64%" highlight style used
#     elsewhere in the same bullet), so plain Find/Replace won't do.

$locate = $d.Content.Duplicate
$locate.Find.Execute(
    "all Black and Asian-American",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$phraseStart = $locate.Start
$phraseEnd = $locate.End

# Overwrite "all Black and Asian-American" (the word " voters" that follows
# stays untouched) with "50M", then apply bold + the same theme color used
# for the other inline-highlighted numbers in this bullet.
$target = $d.Range($phraseStart, $phraseEnd)
$target.Text = "50M"

$highlight = $d.Range($phraseStart, $phraseStart + 3)
$highlight.Font.Bold = $true
$highlight.Font.Color = 5258796
